$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Column B (spending I)
Set-CellText "B2"  "36.24***"
Set-CellText "B3"  "(7.46)"
Set-CellText "B12" "61422"

# Column C (spending II)
Set-CellText "C4"  "0.63***"
Set-CellText "C5"  "(0.12)"
Set-CellText "C12" "62959"

# Column D (spending III)
Set-CellText "D6"  "1.07***"
Set-CellText "D7"  "(0.13)"
Set-CellText "D12" "59057"

# Column E (spending IIII)
Set-CellText "E9"  "(0.40)"
Set-CellText "E12" "60560"

# Column F (spending IIIII)
Set-CellText "F10" "0.16"
Set-CellText "F11" "(0.13)"
Set-CellText "F12" "87775"
